# Updates cryptos list values (Price / Volume(1h) columns) to match the
# latest scrape. Numeric-looking price strings are forced back to text
# (matching the original inline-string cells) via a NumberFormat="@" /
# ClearFormats() round trip, so Excel does not silently coerce them to
# the Number type, while leaving cell formatting untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.720.01"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "2.085.17"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.08"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").Value = "2.394.76"
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("E13").Value = "  +3.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.85"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.799"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.70%  "
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "2.084.87"
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "38.671.94"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.02"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.98"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("E25").Value = "  +2.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.91"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("E27").Value = "  +2.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.138"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.12%  "
$ws.Range("E29").Value = "  +13.10%  "
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +4.34%  "
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.69"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0607"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("E38").Value = "  +2.52%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.85"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("E41").Value = "  +5.59%  "
$ws.Range("D42").Value = "1.542.21"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.34"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.21%  "
$ws.Range("E45").Value = "  +3.66%  "
$ws.Range("E46").Value = "  +8.97%  "
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.13"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  +2.56%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "2.283.11"
$ws.Range("E51").Value = "  +1.88%  "
